$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet had an extra, redundant column (M) in the alcohol measurement
# data; remove it so the data that used to live in column N shifts left
# into column M.
$ws.Columns.Item(13).Delete()

# Reflect the resulting selection on the now-last data column.
$ws.Range("M1").Select()
